$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.59 = 5820.63 pesos`n✅ 5820.63 pesos = 1.59 = 952.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update tasas sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 630
$ws2.Range("O10").Value = 3667
$ws2.Range("N12").Value = 3672
$ws2.Range("O12").Value = 601.01
